$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the diff-comparison header columns: "_old" -> "_FV2310" and
#    "_new" -> "_FV2404" (the merged-AHB columns are being relabelled
#    to the actual format-version identifiers being compared).
# ---------------------------------------------------------------------
$ws.Range("A1").Value()  = "Segmentname_FV2310"
$ws.Range("B1").Value()  = "Segmentgruppe_FV2310"
$ws.Range("C1").Value()  = "Segment_FV2310"
$ws.Range("D1").Value()  = "Datenelement_FV2310"
$ws.Range("E1").Value()  = "Segment ID_FV2310"
$ws.Range("F1").Value()  = "Code_FV2310"
$ws.Range("G1").Value()  = "Qualifier_FV2310"
$ws.Range("H1").Value()  = "Beschreibung_FV2310"
$ws.Range("I1").Value()  = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value()  = "Bedingung_FV2310"
$ws.Range("K1").Value()  = "diff"
$ws.Range("L1").Value()  = "Segmentname_FV2404"
$ws.Range("M1").Value()  = "Segmentgruppe_FV2404"
$ws.Range("N1").Value()  = "Segment_FV2404"
$ws.Range("O1").Value()  = "Datenelement_FV2404"
$ws.Range("P1").Value()  = "Segment ID_FV2404"
$ws.Range("Q1").Value()  = "Code_FV2404"
$ws.Range("R1").Value()  = "Qualifier_FV2404"
$ws.Range("S1").Value()  = "Beschreibung_FV2404"
$ws.Range("T1").Value()  = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value()  = "Bedingung_FV2404"

# ---------------------------------------------------------------------
# 2) Turn the used range A1:U58 into a real Excel Table ("Table1") with
#    an AutoFilter, so the header row can be sorted/filtered.
#
#    The header row already carries bold/fill/border formatting (style
#    index 1); creating the ListObject over a styled header would make
#    Excel bake that look into a one-off conditional-format (dxf) tied
#    to the table's headerRowDxfId. To keep the original shared cell
#    style intact (no extra dxf), the header formatting is stashed on
#    a scratch row, cleared, then restored via a format-only paste once
#    the table exists.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$headerRange.Copy() | Out-Null
$scratchRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats() | Out-Null

$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratchRange.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(100).Delete() | Out-Null

# ---------------------------------------------------------------------
# 3) Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
